$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 1320.5
$ws.Range("I13").Value = 2001
$ws.Range("J13").Value = 1093.6666
$ws.Range("K13").Value = 2001
$ws.Range("L13").Value = 1093.6666
$ws.Range("M13").Value = -1832
$ws.Range("N13").Value = -1431.6666
$ws.Range("H15").Value = 1130.2273
$ws.Range("I15").Value = 1130.2273
$ws.Range("K15").Value = 3390.6819
$ws.Range("M15").Value = -3221.6819
$ws.Range("H17").Value = 1096.6177
$ws.Range("J17").Value = 1121.3636
$ws.Range("L17").Value = 3364.0908
$ws.Range("N17").Value = -3700.0908
$ws.Range("H19").Value = 3644.6
$ws.Range("J19").Value = 3407.3333
$ws.Range("L19").Value = 3407.3333
$ws.Range("N19").Value = -3757.3333
$ws.Range("H41").Value = 16131223
$ws.Range("I41").Value = 508.47058
$ws.Range("K41").Value = 508.47058
$ws.Range("M41").Value = -68.47057999999998
$ws.Range("H53").Value = 17544852
$ws.Range("I53").Value = 37037664
$ws.Range("J53").Value = 1322
$ws.Range("K53").Value = 37037664
$ws.Range("L53").Value = 1322
$ws.Range("M53").Value = -37037027
$ws.Range("N53").Value = -2596
$ws.Range("H62").Value = 11365762
$ws.Range("I62").Value = 13890777
$ws.Range("J62").Value = 3194.5
$ws.Range("K62").Value = 13890777
$ws.Range("L62").Value = 3194.5
$ws.Range("M62").Value = -13890153
$ws.Range("N62").Value = -4442.5
$ws.Range("H64").Value = 8668.166999999999
$ws.Range("I64").Value = 5003
$ws.Range("J64").Value = 10500.75
$ws.Range("K64").Value = 5003
$ws.Range("L64").Value = 10500.75
$ws.Range("M64").Value = -4755
$ws.Range("N64").Value = -10996.75
$ws.Range("H65").Value = 11365762
$ws.Range("I65").Value = 13890777
$ws.Range("J65").Value = 3194.5
$ws.Range("K65").Value = 69453885
$ws.Range("L65").Value = 15972.5
$ws.Range("M65").Value = -69450765
$ws.Range("N65").Value = -22212.5
$ws.Range("H67").Value = 8668.166999999999
$ws.Range("I67").Value = 5003
$ws.Range("J67").Value = 10500.75
$ws.Range("K67").Value = 5003
$ws.Range("L67").Value = 10500.75
$ws.Range("M67").Value = -4145
$ws.Range("N67").Value = -12216.75
$ws.Range("H74").Value = 12125.786
$ws.Range("I74").Value = 11135.462
$ws.Range("K74").Value = 11135.462
$ws.Range("M74").Value = -10199.462
$ws.Range("H77").Value = 12125.786
$ws.Range("I77").Value = 11135.462
$ws.Range("K77").Value = 55677.31
$ws.Range("M77").Value = -50997.31
$ws.Range("H103").Value = 239
$ws.Range("I103").Value = 186.8
$ws.Range("K103").Value = 560.4000000000001
$ws.Range("M103").Value = 25.59999999999991
$ws.Range("H107").Value = 85024.5
$ws.Range("I107").Value = 112710.445
$ws.Range("J107").Value = 1966.6666
$ws.Range("K107").Value = 112710.445
$ws.Range("L107").Value = 1966.6666
$ws.Range("M107").Value = -110790.445
$ws.Range("N107").Value = -5806.6666
$ws.Range("H112").Value = 2462.0454
$ws.Range("J112").Value = 2803.611
$ws.Range("L112").Value = 8410.832999999999
$ws.Range("N112").Value = -10626.833
$ws.Range("H116").Value = 4399.8
$ws.Range("I116").Value = 4249.75
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 4249.75
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -807.75
$ws.Range("N116").Value = -11884
$ws.Range("H128").Value = 87358.87
$ws.Range("J128").Value = 87358.87
$ws.Range("L128").Value = 87358.87
$ws.Range("N128").Value = -97318.87
$ws.Range("H132").Value = 4803.675
$ws.Range("I132").Value = 4936
$ws.Range("K132").Value = 14808
$ws.Range("M132").Value = -12278
$ws.Range("H133").Value = 75005
$ws.Range("J133").Value = 75005
$ws.Range("L133").Value = 75005
$ws.Range("N133").Value = -85125
$ws.Range("H134").Value = 71052.52
$ws.Range("J134").Value = 71052.52
$ws.Range("L134").Value = 71052.52
$ws.Range("N134").Value = -81192.52
$ws.Range("H135").Value = 3000.7
$ws.Range("I135").Value = 1358.1428
$ws.Range("K135").Value = 12223.2852
$ws.Range("M135").Value = -9688.2852
$ws.Range("H138").Value = 4159.985
$ws.Range("J138").Value = 5012.84
$ws.Range("L138").Value = 15038.52
$ws.Range("N138").Value = -25318.52
$ws.Range("H139").Value = 115000
$ws.Range("J139").Value = 115000
$ws.Range("L139").Value = 115000
$ws.Range("N139").Value = -125280
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3688.4792
$ws.Range("I32").Value = 3554.1914
$ws.Range("K32").Value = 3554.1914
$ws.Range("M32").Value = -3267.1914
$ws.Range("H61").Value = 4301.3184
$ws.Range("I61").Value = 1804.5385
$ws.Range("J61").Value = 7907.778
$ws.Range("K61").Value = 1804.5385
$ws.Range("L61").Value = 7907.778
$ws.Range("M61").Value = -1592.5385
$ws.Range("N61").Value = -8331.778
$ws.Range("H122").Value = 7046
$ws.Range("I122").Value = 8950
$ws.Range("J122").Value = 6502
$ws.Range("K122").Value = 26850
$ws.Range("L122").Value = 19506
$ws.Range("M122").Value = -24400
$ws.Range("N122").Value = -24406
$ws.Range("H132").Value = 5194.222
$ws.Range("I132").Value = 2698.4119
$ws.Range("K132").Value = 8095.2357
$ws.Range("M132").Value = -5565.2357
$ws.Range("H136").Value = 4301.3184
$ws.Range("I136").Value = 1804.5385
$ws.Range("J136").Value = 7907.778
$ws.Range("K136").Value = 5413.6155
$ws.Range("L136").Value = 23723.334
$ws.Range("M136").Value = -2863.6155
$ws.Range("N136").Value = -28823.334

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1227.2903
$ws.Range("I20").Value = 871.6
$ws.Range("J20").Value = 1874
$ws.Range("K20").Value = 871.6
$ws.Range("L20").Value = 1874
$ws.Range("M20").Value = -624.6
$ws.Range("N20").Value = -2368
$ws.Range("H80").Value = 2451.1428
$ws.Range("J80").Value = 1832.2
$ws.Range("L80").Value = 1832.2
$ws.Range("N80").Value = -3828.2
$ws.Range("H82").Value = 48453.125
$ws.Range("I82").Value = 22465
$ws.Range("J82").Value = 91766.664
$ws.Range("K82").Value = 22465
$ws.Range("L82").Value = 91766.664
$ws.Range("M82").Value = -22082
$ws.Range("N82").Value = -92532.664
$ws.Range("H83").Value = 2451.1428
$ws.Range("J83").Value = 1832.2
$ws.Range("L83").Value = 9161
$ws.Range("N83").Value = -19145
$ws.Range("H85").Value = 48453.125
$ws.Range("I85").Value = 22465
$ws.Range("J85").Value = 91766.664
$ws.Range("K85").Value = 22465
$ws.Range("L85").Value = 91766.664
$ws.Range("M85").Value = -21139
$ws.Range("N85").Value = -94418.664
$ws.Range("H99").Value = 2683
$ws.Range("I99").Value = 2360.25
$ws.Range("K99").Value = 2360.25
$ws.Range("M99").Value = -862.25
$ws.Range("H132").Value = 66666.664
$ws.Range("J132").Value = 66666.664
$ws.Range("L132").Value = 66666.664
$ws.Range("N132").Value = -76786.664
$ws.Range("H134").Value = 3510.5334
$ws.Range("I134").Value = 2474.125
$ws.Range("K134").Value = 7422.375
$ws.Range("M134").Value = -4887.375

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3465.5
$ws.Range("I31").Value = 2817.5557
$ws.Range("J31").Value = 3914.077
$ws.Range("K31").Value = 2817.5557
$ws.Range("L31").Value = 3914.077
$ws.Range("M31").Value = -2522.5557
$ws.Range("N31").Value = -4504.077
$ws.Range("H34").Value = 3465.5
$ws.Range("I34").Value = 2817.5557
$ws.Range("J34").Value = 3914.077
$ws.Range("K34").Value = 2817.5557
$ws.Range("L34").Value = 3914.077
$ws.Range("M34").Value = -2615.5557
$ws.Range("N34").Value = -4318.077
$ws.Range("H58").Value = 349744.12
$ws.Range("J58").Value = 5833.3335
$ws.Range("L58").Value = 5833.3335
$ws.Range("N58").Value = -6239.3335
$ws.Range("H99").Value = 4991.7334
$ws.Range("I99").Value = 3619.8
$ws.Range("J99").Value = 5677.7
$ws.Range("K99").Value = 3619.8
$ws.Range("L99").Value = 5677.7
$ws.Range("M99").Value = -2121.8
$ws.Range("N99").Value = -8673.700000000001
$ws.Range("H105").Value = 1223.52
$ws.Range("J105").Value = 1429
$ws.Range("L105").Value = 1429
$ws.Range("N105").Value = -4923
$ws.Range("H122").Value = 2380.682
$ws.Range("I122").Value = 1362.4615
$ws.Range("J122").Value = 3851.4443
$ws.Range("K122").Value = 4087.3845
$ws.Range("L122").Value = 11554.3329
$ws.Range("M122").Value = -1637.3845
$ws.Range("N122").Value = -16454.3329
$ws.Range("H126").Value = 4991.7334
$ws.Range("I126").Value = 3619.8
$ws.Range("J126").Value = 5677.7
$ws.Range("K126").Value = 10859.4
$ws.Range("L126").Value = 17033.1
$ws.Range("M126").Value = -8389.400000000001
$ws.Range("N126").Value = -21973.1
$ws.Range("H132").Value = 3799.7878
$ws.Range("I132").Value = 2834.739
$ws.Range("K132").Value = 8504.217000000001
$ws.Range("M132").Value = -5974.217000000001
$ws.Range("H134").Value = 4271.778
$ws.Range("I134").Value = 3467.3928
$ws.Range("K134").Value = 10402.1784
$ws.Range("M134").Value = -7867.178400000001
$ws.Range("H136").Value = 349744.12
$ws.Range("J136").Value = 5833.3335
$ws.Range("L136").Value = 17500.0005
$ws.Range("N136").Value = -22600.0005

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1147.2
$ws.Range("I3").Value = 1147.2
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3441.6
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -3329.6
$ws.Range("N3").ClearContents()
$ws.Range("H4").Value = 3041815.8
$ws.Range("I4").Value = 1851330.2
$ws.Range("K4").Value = 5553990.6
$ws.Range("M4").Value = -5553878.6
$ws.Range("H132").Value = 2772.65
$ws.Range("I132").Value = 1100
$ws.Range("K132").Value = 9900
$ws.Range("M132").Value = -7370
$ws.Range("H140").Value = 4069.7144
$ws.Range("I140").Value = 4069.7144
$ws.Range("K140").Value = 12209.1432
$ws.Range("M140").Value = -7029.143199999999

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 550000000
$ws.Range("I21").Value = 1000000000
$ws.Range("J21").Value = 100000000
$ws.Range("K21").Value = 1000000000
$ws.Range("L21").Value = 100000000
$ws.Range("M21").Value = -999999827
$ws.Range("N21").Value = -100000346
$ws.Range("H24").Value = 2527500
$ws.Range("J24").Value = 3350000
$ws.Range("L24").Value = 3350000
$ws.Range("N24").Value = -3350346
$ws.Range("H30").Value = 550000000
$ws.Range("I30").Value = 1000000000
$ws.Range("J30").Value = 100000000
$ws.Range("K30").Value = 1000000000
$ws.Range("L30").Value = 100000000
$ws.Range("M30").Value = -999999895
$ws.Range("N30").Value = -100000210
$ws.Range("H31").Value = 2949.6667
$ws.Range("I31").Value = 1925
$ws.Range("J31").Value = 4999
$ws.Range("K31").Value = 1925
$ws.Range("L31").Value = 4999
$ws.Range("M31").Value = -1633
$ws.Range("N31").Value = -5583
$ws.Range("H37").Value = 2949.6667
$ws.Range("I37").Value = 1925
$ws.Range("J37").Value = 4999
$ws.Range("K37").Value = 1925
$ws.Range("L37").Value = 4999
$ws.Range("M37").Value = -1648
$ws.Range("N37").Value = -5553
$ws.Range("H80").Value = 3006802.2
$ws.Range("I80").Value = 4005602.5
$ws.Range("J80").Value = 2008002.2
$ws.Range("K80").Value = 4005602.5
$ws.Range("L80").Value = 2008002.2
$ws.Range("M80").Value = -4004604.5
$ws.Range("N80").Value = -2009998.2
$ws.Range("H83").Value = 3006802.2
$ws.Range("I83").Value = 4005602.5
$ws.Range("J83").Value = 2008002.2
$ws.Range("K83").Value = 20028012.5
$ws.Range("L83").Value = 10040011
$ws.Range("M83").Value = -20023020.5
$ws.Range("N83").Value = -10049995
$ws.Range("H111").Value = 37824.668
$ws.Range("J111").Value = 37824.668
$ws.Range("L111").Value = 37824.668
$ws.Range("N111").Value = -43958.668
$ws.Range("H132").Value = 241122.1
$ws.Range("I132").Value = 280314.7
$ws.Range("K132").Value = 840944.1000000001
$ws.Range("M132").Value = -838414.1000000001

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 50000
$ws.Range("J20").Value = 50000
$ws.Range("L20").Value = 50000
$ws.Range("N20").Value = -50452
$ws.Range("H22").Value = 1206.5769
$ws.Range("I22").Value = 681.7273
$ws.Range("J22").Value = 1591.4667
$ws.Range("K22").Value = 681.7273
$ws.Range("L22").Value = 1591.4667
$ws.Range("M22").Value = -386.7273
$ws.Range("N22").Value = -2181.4667
$ws.Range("H23").Value = 67775
$ws.Range("I23").Value = 111100
$ws.Range("K23").Value = 111100
$ws.Range("M23").Value = -110870
$ws.Range("H27").Value = 1206.5769
$ws.Range("I27").Value = 681.7273
$ws.Range("J27").Value = 1591.4667
$ws.Range("K27").Value = 681.7273
$ws.Range("L27").Value = 1591.4667
$ws.Range("M27").Value = -574.7273
$ws.Range("N27").Value = -1805.4667
$ws.Range("H46").Value = 5068.316
$ws.Range("I46").Value = 3733.1667
$ws.Range("K46").Value = 3733.1667
$ws.Range("M46").Value = -3545.1667
$ws.Range("H55").Value = 540.3415
$ws.Range("I55").Value = 449.625
$ws.Range("J55").Value = 862.8889
$ws.Range("K55").Value = 449.625
$ws.Range("L55").Value = 862.8889
$ws.Range("M55").Value = -276.625
$ws.Range("N55").Value = -1208.8889
$ws.Range("H68").Value = 156120.42
$ws.Range("I68").Value = 18170.8
$ws.Range("J68").Value = 500994.5
$ws.Range("K68").Value = 18170.8
$ws.Range("L68").Value = 500994.5
$ws.Range("M68").Value = -17421.8
$ws.Range("N68").Value = -502492.5
$ws.Range("H71").Value = 156120.42
$ws.Range("I71").Value = 18170.8
$ws.Range("J71").Value = 500994.5
$ws.Range("K71").Value = 90854
$ws.Range("L71").Value = 2504972.5
$ws.Range("M71").Value = -87110
$ws.Range("N71").Value = -2512460.5
$ws.Range("H82").Value = 2324.4666
$ws.Range("I82").Value = 2376.375
$ws.Range("K82").Value = 2376.375
$ws.Range("M82").Value = -2015.375
$ws.Range("H85").Value = 2324.4666
$ws.Range("I85").Value = 2376.375
$ws.Range("K85").Value = 2376.375
$ws.Range("M85").Value = -1128.375
$ws.Range("H132").Value = 7999.375
$ws.Range("I132").Value = 4972.5
$ws.Range("J132").Value = 9008.333000000001
$ws.Range("K132").Value = 14917.5
$ws.Range("L132").Value = 27024.999
$ws.Range("M132").Value = -12387.5
$ws.Range("N132").Value = -32084.999
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 4283
$ws.Range("I136").Value = 3633.182
$ws.Range("K136").Value = 10899.546
$ws.Range("M136").Value = -8349.545999999998

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 29999
$ws.Range("I51").Value = 29999
$ws.Range("K51").Value = 29999
$ws.Range("M51").Value = -29489
$ws.Range("H81").Value = 57461
$ws.Range("I81").Value = 104973.5
$ws.Range("J81").Value = 9948.5
$ws.Range("K81").Value = 209947
$ws.Range("L81").Value = 19897
$ws.Range("M81").Value = -208886
$ws.Range("N81").Value = -22019
$ws.Range("H84").Value = 57461
$ws.Range("I84").Value = 104973.5
$ws.Range("J84").Value = 9948.5
$ws.Range("K84").Value = 1049735
$ws.Range("L84").Value = 99485
$ws.Range("M84").Value = -1044431
$ws.Range("N84").Value = -110093
$ws.Range("H107").Value = 39152.223
$ws.Range("I107").Value = 50017
$ws.Range("J107").Value = 1125.5
$ws.Range("K107").Value = 150051
$ws.Range("L107").Value = 3376.5
$ws.Range("M107").Value = -148131
$ws.Range("N107").Value = -7216.5
$ws.Range("H132").Value = 4092.5
$ws.Range("I132").Value = 3756.4707
$ws.Range("K132").Value = 11269.4121
$ws.Range("M132").Value = -8739.4121
